$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 623
$wsExhibit.Range("F4").Value = 619
$wsExhibit.Range("F5").Value = 546
$wsExhibit.Range("F6").Value = 304
$wsExhibit.Range("F7").Value = 2742
$wsExhibit.Range("F8").Value = 464
$wsExhibit.Range("F9").Value = 7657
$wsExhibit.Range("F11").Value = 468
$wsExhibit.Range("F12").Value = 34
$wsExhibit.Range("F13").Value = 300
$wsExhibit.Range("F14").Value = 43

# Sheet "全部类型" (sheet4): update "想去人数" (column F) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 623
$wsAll.Range("F4").Value = 619
$wsAll.Range("F5").Value = 546
$wsAll.Range("F6").Value = 304
$wsAll.Range("F9").Value = 2742
$wsAll.Range("F10").Value = 464
$wsAll.Range("F11").Value = 7657
$wsAll.Range("F13").Value = 468
$wsAll.Range("F14").Value = 34
$wsAll.Range("F17").Value = 300
$wsAll.Range("F18").Value = 43

$wb.Save()
